$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column G entirely (values + shift dimension from A1:G26 to A1:F26)
$ws.Range("G1").EntireColumn.Delete()

# Update changed values per row (columns B:F)
$ws.Range("B2").Value = 0.051
$ws.Range("C2").Value = 0.15
$ws.Range("D2").Value = 0.376
$ws.Range("E2").Value = 0.74
$ws.Range("F2").Value = 0.991

$ws.Range("B3").Value = 0.048
$ws.Range("C3").Value = 0.143
$ws.Range("D3").Value = 0.361
$ws.Range("E3").Value = 0.731
$ws.Range("F3").Value = 0.989

$ws.Range("B4").Value = 0.048
$ws.Range("C4").Value = 0.142
$ws.Range("D4").Value = 0.366
$ws.Range("E4").Value = 0.732
$ws.Range("F4").Value = 0.99

$ws.Range("B5").Value = 0.048
$ws.Range("C5").Value = 0.144
$ws.Range("D5").Value = 0.369
$ws.Range("E5").Value = 0.734
$ws.Range("F5").Value = 0.991

$ws.Range("B6").Value = 0.052
$ws.Range("C6").Value = 0.15
$ws.Range("D6").Value = 0.377
$ws.Range("E6").Value = 0.74
$ws.Range("F6").Value = 0.991

$ws.Range("B7").Value = 0.048
$ws.Range("C7").Value = 0.142
$ws.Range("D7").Value = 0.366
$ws.Range("E7").Value = 0.732
$ws.Range("F7").Value = 0.99

$ws.Range("B8").Value = 0.067
$ws.Range("C8").Value = 0.163
$ws.Range("D8").Value = 0.381
$ws.Range("E8").Value = 0.744
$ws.Range("F8").Value = 0.993

$ws.Range("B9").Value = 0.072
$ws.Range("C9").Value = 0.17
$ws.Range("D9").Value = 0.39
$ws.Range("E9").Value = 0.755
$ws.Range("F9").Value = 0.993

$ws.Range("B10").Value = 0.1
$ws.Range("C10").Value = 0.211
$ws.Range("D10").Value = 0.548
$ws.Range("E10").Value = 0.985
$ws.Range("F10").Value = 0.993

$ws.Range("B11").Value = 0.049
$ws.Range("C11").Value = 0.144
$ws.Range("D11").Value = 0.377
$ws.Range("E11").Value = 0.751
$ws.Range("F11").Value = 0.991

$ws.Range("B12").Value = 0.049
$ws.Range("C12").Value = 0.144
$ws.Range("D12").Value = 0.377
$ws.Range("E12").Value = 0.751
$ws.Range("F12").Value = 0.991

$ws.Range("B13").Value = 0.048
$ws.Range("C13").Value = 0.142
$ws.Range("D13").Value = 0.368
$ws.Range("E13").Value = 0.734
$ws.Range("F13").Value = 0.991

$ws.Range("B14").Value = 0.127
$ws.Range("C14").Value = 0.233
$ws.Range("D14").Value = 0.583
$ws.Range("E14").Value = 0.99

$ws.Range("B15").Value = 0.048
$ws.Range("C15").Value = 0.144
$ws.Range("D15").Value = 0.373
$ws.Range("E15").Value = 0.745
$ws.Range("F15").Value = 0.991

$ws.Range("B16").Value = 0.048
$ws.Range("C16").Value = 0.142
$ws.Range("D16").Value = 0.366
$ws.Range("E16").Value = 0.732
$ws.Range("F16").Value = 0.99

$ws.Range("B17").Value = 0.048
$ws.Range("C17").Value = 0.142
$ws.Range("D17").Value = 0.366
$ws.Range("E17").Value = 0.732
$ws.Range("F17").Value = 0.99

$ws.Range("B18").Value = 0.13
$ws.Range("C18").Value = 0.366

$ws.Range("B19").Value = 0.144

$ws.Range("B20").Value = 0.048
$ws.Range("C20").Value = 0.142
$ws.Range("D20").Value = 0.366
$ws.Range("E20").Value = 0.737
$ws.Range("F20").Value = 0.994

$ws.Range("B21").Value = 0.107
$ws.Range("C21").Value = 0.228
$ws.Range("D21").Value = 0.732

$ws.Range("B22").Value = 0.048
$ws.Range("C22").Value = 0.142
$ws.Range("D22").Value = 0.366
$ws.Range("E22").Value = 0.732
$ws.Range("F22").Value = 0.99

$ws.Range("B23").Value = 0.052
$ws.Range("C23").Value = 0.151
$ws.Range("D23").Value = 0.392
$ws.Range("E23").Value = 0.824
$ws.Range("F23").Value = 1

$ws.Range("B24").Value = 0.054
$ws.Range("C24").Value = 0.165
$ws.Range("D24").Value = 0.469
$ws.Range("E24").Value = 0.804
$ws.Range("F24").Value = 1

$ws.Range("B25").Value = 0.048
$ws.Range("C25").Value = 0.14
$ws.Range("D25").Value = 0.323
$ws.Range("E25").Value = 0.742
$ws.Range("F25").Value = 0.988

$ws.Range("B26").Value = 0.052
$ws.Range("C26").Value = 0.158
$ws.Range("D26").Value = 0.407
$ws.Range("E26").Value = 0.785
$ws.Range("F26").Value = 1
